$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain literal
# text (matching the source inlineStr cells) - force Text format first so
# Excel does not auto-convert them (e.g. "1.00" -> 1, "0.733" -> 0.733 losing the string type).
$textCells = @('D4', 'D5', 'D6', 'D9', 'D12', 'D13', 'D16', 'D18', 'D21', 'D22', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.357.32'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '3.912.67'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '485.79'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '145.68'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.733'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').Value = '43.09'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '10.73'
$ws.Range('E13').Value = '  +2.76%  '
$ws.Range('D14').Value = '4.536.94'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').Value = '3.931.29'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '14.33'
$ws.Range('E16').Value = '  -5.43%  '
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '19.99'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = '68.396.58'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('D21').Value = '432.19'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').Value = '15.17'
$ws.Range('E22').Value = '  +4.34%  '
$ws.Range('E23').Value = '  +2.17%  '
$ws.Range('D24').Value = '88.04'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = '11.48'
$ws.Range('E25').Value = '  +16.44%  '
$ws.Range('D26').Value = '11.26'
$ws.Range('E26').Value = '  +11.91%  '
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').Value = '37.86'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').Value = '5.69'
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').Value = '718.24'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('E31').Value = '  +3.11%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('E33').Value = '  +4.21%  '
$ws.Range('D34').Value = '6.16'
$ws.Range('E34').Value = '  +14.13%  '
$ws.Range('D35').Value = '41.37'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').Value = '0.0₃0876'
$ws.Range('E36').Value = '  +4.24%  '
$ws.Range('D37').Value = '60.97'
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('D38').Value = '0.147'
$ws.Range('E38').Value = '  -4.12%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0503'
$ws.Range('E39').Value = '  +5.72%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '3.02'
$ws.Range('E41').Value = '  +19.10%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '0.392'
$ws.Range('E42').Value = '  +15.33%  '
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('E44').Value = '  +5.11%  '
$ws.Range('E45').Value = '  +5.27%  '
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('E49').Value = '  -4.93%  '
$ws.Range('D50').Value = '144.95'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('D51').Value = '0.0₆0337'
$ws.Range('E51').Value = '  +23.59%  '
